# Applies the betexplorer re-scrape update for serbia/prva-liga 2023-2024:
#  - several existing matchdays had their row order reshuffled (values for
#    columns F:V swap between rows that share the same match date in A:E)
#  - three new matches (rows 111-113) are appended at the end
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: reorder existing rows within their matchday groups ---
# (columns A:E -- index/country/tournament/season/date -- are untouched;
#  only F:V, the match details, move between rows)

# g1 (25-27, 2023-08-26 round)
$v25 = $ws.Range("F25:V25").Value2
$v26 = $ws.Range("F26:V26").Value2
$v27 = $ws.Range("F27:V27").Value2
$ws.Range("F25:V25").Value = $v26
$ws.Range("F26:V26").Value = $v27
$ws.Range("F27:V27").Value = $v25

# g2 (29 & 31, 2023-08-27 round)
$v29 = $ws.Range("F29:V29").Value2
$v31 = $ws.Range("F31:V31").Value2
$ws.Range("F29:V29").Value = $v31
$ws.Range("F31:V31").Value = $v29

# g3 (65-67, 2023-10-01 round)
$v65 = $ws.Range("F65:V65").Value2
$v66 = $ws.Range("F66:V66").Value2
$v67 = $ws.Range("F67:V67").Value2
$ws.Range("F65:V65").Value = $v66
$ws.Range("F66:V66").Value = $v67
$ws.Range("F67:V67").Value = $v65

# g4 (79-80, 2023-10-14 round)
$v79 = $ws.Range("F79:V79").Value2
$v80 = $ws.Range("F80:V80").Value2
$ws.Range("F79:V79").Value = $v80
$ws.Range("F80:V80").Value = $v79

# g5 (82-83, 2023-10-15 round)
$v82 = $ws.Range("F82:V82").Value2
$v83 = $ws.Range("F83:V83").Value2
$ws.Range("F82:V82").Value = $v83
$ws.Range("F83:V83").Value = $v82

# g6 (84 & 86, 2023-10-16 round)
$v84 = $ws.Range("F84:V84").Value2
$v86 = $ws.Range("F86:V86").Value2
$ws.Range("F84:V84").Value = $v86
$ws.Range("F86:V86").Value = $v84

# g7 (97-101, 2023-10-28 round)
$v97 = $ws.Range("F97:V97").Value2
$v98 = $ws.Range("F98:V98").Value2
$v99 = $ws.Range("F99:V99").Value2
$v100 = $ws.Range("F100:V100").Value2
$v101 = $ws.Range("F101:V101").Value2
$ws.Range("F97:V97").Value = $v98
$ws.Range("F98:V98").Value = $v101
$ws.Range("F99:V99").Value = $v97
$ws.Range("F100:V100").Value = $v99
$ws.Range("F101:V101").Value = $v100

# g8 (108-109, 2023-11-06 round)
$v108 = $ws.Range("F108:V108").Value2
$v109 = $ws.Range("F109:V109").Value2
$ws.Range("F108:V108").Value = $v109
$ws.Range("F109:V109").Value = $v108

# --- Step 2: append three new match rows at the bottom (111-113) ---
# Copy formatting (styles) from the last existing row so the new rows
# get the same number formats / fonts / borders as the rest of the table.
$ws.Range("A110:V110").Copy()
$ws.Range("A111:V113").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# new row 111
$ws.Range("A111").Value = 110
$ws.Range("B111").Value = "serbia"
$ws.Range("C111").Value = "prva-liga"
$ws.Range("D111").Value = "2023-2024"
$ws.Range("E111").Value = 45241.54166666666
$ws.Range("F111").Value = "Radnicki S. Mitrovica"
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = "Mladost GAT"
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 2.56
$ws.Range("K111").Value = "11/11/2023 02:13"
$ws.Range("L111").Value = 2.04
$ws.Range("M111").Value = "11/11/2023 12:53"
$ws.Range("N111").Value = 2.74
$ws.Range("O111").Value = "11/11/2023 02:13"
$ws.Range("P111").Value = 2.73
$ws.Range("Q111").Value = "11/11/2023 12:53"
$ws.Range("R111").Value = 2.81
$ws.Range("S111").Value = "11/11/2023 02:13"
$ws.Range("T111").Value = 4.03
$ws.Range("U111").Value = "11/11/2023 12:53"
$ws.Range("V111").Value = "https://www.betexplorer.com/football/serbia/prva-liga/radnicki-s-mitrovica-mladost-gat/88pX5ob4/"

# new row 112
$ws.Range("A112").Value = 111
$ws.Range("B112").Value = "serbia"
$ws.Range("C112").Value = "prva-liga"
$ws.Range("D112").Value = "2023-2024"
$ws.Range("E112").Value = 45241.54166666666
$ws.Range("F112").Value = "FK Indjija"
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = "Smederevo"
$ws.Range("I112").Value = 2
$ws.Range("J112").Value = 1.88
$ws.Range("K112").Value = "11/11/2023 02:13"
$ws.Range("L112").Value = 1.99
$ws.Range("M112").Value = "11/11/2023 12:48"
$ws.Range("N112").Value = 2.98
$ws.Range("O112").Value = "11/11/2023 02:13"
$ws.Range("P112").Value = 2.8
$ws.Range("Q112").Value = "11/11/2023 12:48"
$ws.Range("R112").Value = 3.91
$ws.Range("S112").Value = "11/11/2023 02:13"
$ws.Range("T112").Value = 4.05
$ws.Range("U112").Value = "11/11/2023 12:35"
$ws.Range("V112").Value = "https://www.betexplorer.com/football/serbia/prva-liga/indjija-smederevo/UHoT6Rrb/"

# new row 113
$ws.Range("A113").Value = 112
$ws.Range("B113").Value = "serbia"
$ws.Range("C113").Value = "prva-liga"
$ws.Range("D113").Value = "2023-2024"
$ws.Range("E113").Value = 45241.54166666666
$ws.Range("F113").Value = "Tekstilac Odzaci"
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = "Macva"
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1.8
$ws.Range("K113").Value = "11/11/2023 02:13"
$ws.Range("L113").Value = 2.15
$ws.Range("M113").Value = "11/11/2023 12:58"
$ws.Range("N113").Value = 2.98
$ws.Range("O113").Value = "11/11/2023 02:13"
$ws.Range("P113").Value = 2.41
$ws.Range("Q113").Value = "11/11/2023 12:58"
$ws.Range("R113").Value = 4.32
$ws.Range("S113").Value = "11/11/2023 02:13"
$ws.Range("T113").Value = 4.44
$ws.Range("U113").Value = "11/11/2023 12:58"
$ws.Range("V113").Value = "https://www.betexplorer.com/football/serbia/prva-liga/tekstilac-odzaci-macva-sabac/2Vqy55DA/"
